$d = $word.ActiveDocument

# The paragraph ending in "unlike other cities." also carries the _GoBack
# bookmark at its very end. New text must land *before* that bookmark, but
# this engine's InsertAfter() always appends past a trailing bookmark. So:
# remove the bookmark first, append the new italic sentence (as the four
# separate runs the diff specifies), then re-create the _GoBack bookmark
# inside the paragraph that used to be the lone empty <w:p/> right after it.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "unlike other cities."
$find.Execute() | Out-Null
$r = $find.Parent
$targetPara = $r.Paragraphs(1)
$r.Collapse(0)

$newRuns = @(
    " In addition, I will keep the record of Gillette city and Rock Springs city because",
    " t",
    "he dataset is small and the cities have",
    " an outlier in only one field."
)

foreach ($t in $newRuns) {
    $startPos = $r.Start
    $r.InsertAfter($t)
    $newRange = $d.Range($startPos, $startPos + $t.Length)
    $newRange.Font.Italic = $true
    $r = $d.Range($startPos + $t.Length, $startPos + $t.Length)
}

# Re-home the _GoBack bookmark in the next paragraph (formerly an empty
# <w:p/>). A bookmark can't be Add()-ed directly at a collapsed position
# that sits right at a paragraph's end, so insert a throwaway character,
# wrap the bookmark around it, then delete the character again - this
# leaves the bookmark collapsed in place, same as the original markup.
$emptyPara = $targetPara.Next()
$pr = $emptyPara.Range
$pr.Collapse(1)
$pr.InsertAfter("Z")

$emptyPara = $targetPara.Next()
$placeholder = $d.Range($emptyPara.Range.Start, $emptyPara.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $placeholder)

$placeholder = $d.Range($emptyPara.Range.Start, $emptyPara.Range.Start + 1)
$placeholder.Text = ""
